$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (l0_neurons = 0, 4 neurons)
$ws.Range("C2").Value = 0.9852
$ws.Range("D2").Value = 0.9851713305628002
$ws.Range("E2").Value = 0.9960079840319361
$ws.Range("F2").Value = 0.9899598393574297
$ws.Range("G2").Value = 0.973630831643002
$ws.Range("H2").Value = 0.9771598808341607
$ws.Range("I2").Value = 0.9890981169474727

# Row 3 (l0_neurons = 1, 16 neurons)
$ws.Range("C3").Value = 0.9944
$ws.Range("D3").Value = 0.9943959943959945
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.9990009990009989
$ws.Range("G3").Value = 0.9869869869869869
$ws.Range("H3").Value = 0.988988988988989
$ws.Range("I3").Value = 0.997002997002997

# Row 4 (l0_neurons = 2, 32 neurons)
$ws.Range("C4").Value = 0.9948
$ws.Range("D4").Value = 0.9947983703597827
$ws.Range("G4").Value = 0.9880239520958083
$ws.Range("H4").Value = 0.9889669007021062

# Row 5 (l0_neurons = 3, 64 neurons)
$ws.Range("D5").Value = 0.9939927843783458
$ws.Range("E5").Value = 1
$ws.Range("G5").Value = 0.9859719438877755
$ws.Range("H5").Value = 0.9869869869869869
$ws.Range("I5").Value = 0.998003992015968

# Row 6 (l0_neurons = 4, 56 neurons)
$ws.Range("C6").Value = 0.9948
$ws.Range("D6").Value = 0.9947983703597827
$ws.Range("G6").Value = 0.9880239520958083
$ws.Range("H6").Value = 0.9889669007021062
$ws.Range("I6").Value = 0.998
